$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant (copies only the cell formatting, reusing existing
# style entries instead of minting new ones).
$xlPasteFormats = -4122

# Target values for the new "2020" column (L), each mirroring the number
# format / style already used by the corresponding "2019" column (K) cell.
$values = @{
    4  = 2020
    5  = 1.2
    6  = 1.7
    7  = 0.4
    8  = 3.3
    9  = 3.9
    10 = 2.4
    11 = 95.5
    12 = 94.4
    13 = 97.2
}

foreach ($row in 4..13) {
    $src = $ws.Range("K$row")
    $dst = $ws.Range("L$row")

    # Match K's formatting (number format, borders, font) on L before writing
    # the value, so the new cell renders identically to its neighbour.
    $src.Copy()
    $dst.PasteSpecial($xlPasteFormats)

    $dst.Value = $values[$row]
}

$excel.CutCopyMode = $false

# Scroll the view right and select the newly-populated column so the sheet
# opens focused on the data that changed.
$ws.Range("L4:L13").Select()
